$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComponentsAndLicenses")

# Deprecate license guessing (stage 2): drop the "Guessed License URL" /
# "Guessed License URL AuditInfo" columns (N:O) from the
# ComponentsAndLicenses sheet.
$ws.Range("N1:O1").EntireColumn.Delete()

# Reflect the selection left behind where the deleted columns used to be
# (matches the post-delete selection Excel leaves active).
$ws.Activate()
[void]$ws.Range("N1:O1048576").Select()
